# Update the "想去人数" (interest count) figures on the 展览 and 全部类型
# sheets to reflect the newly generated output.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6392
$ws1.Range("F8").Value = 24
$ws1.Range("F9").Value = 75
$ws1.Range("F15").Value = 3110
$ws1.Range("F18").Value = 1775
$ws1.Range("F19").Value = 22

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6392
$ws4.Range("F9").Value = 24
$ws4.Range("F10").Value = 75
$ws4.Range("F16").Value = 3110
$ws4.Range("F19").Value = 1775
$ws4.Range("F20").Value = 22
